$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,2).Value = 1.02
$ws.Cells.Item(2,3).Value = 1.028535316525087
$ws.Cells.Item(2,4).Value = 1.031602376488455
$ws.Cells.Item(2,5).Value = 1.02840618210332
$ws.Cells.Item(2,6).Value = 1.036441763206023
$ws.Cells.Item(2,9).Value = 1.032628309558649
$ws.Cells.Item(2,10).Value = 1.033687095665741
$ws.Cells.Item(2,11).Value = 1.0344100198557
$ws.Cells.Item(2,12).Value = 1.031223091868364
$ws.Cells.Item(2,13).Value = 1.039235494404243
$ws.Cells.Item(2,14).Value = 1.035155049907163

$ws.Cells.Item(3,2).Value = 1.02
$ws.Cells.Item(3,3).Value = 1.03002376189744
$ws.Cells.Item(3,4).Value = 1.032711804077479
$ws.Cells.Item(3,5).Value = 1.029691197379272
$ws.Cells.Item(3,6).Value = 1.038819834016068
$ws.Cells.Item(3,9).Value = 1.033066129558511
$ws.Cells.Item(3,10).Value = 1.034813430928226
$ws.Cells.Item(3,11).Value = 1.035327239432978
$ws.Cells.Item(3,12).Value = 1.032314749258906
$ws.Cells.Item(3,13).Value = 1.041419012674996
$ws.Cells.Item(3,14).Value = 1.036282984694904

$ws.Cells.Item(4,2).Value = 1.02
$ws.Cells.Item(4,3).Value = 1.030984754399783
$ws.Cells.Item(4,4).Value = 1.033427570199331
$ws.Cells.Item(4,5).Value = 1.03052095537228
$ws.Cells.Item(4,6).Value = 1.040351132919755
$ws.Cells.Item(4,9).Value = 1.033346414036041
$ws.Cells.Item(4,10).Value = 1.035539768421997
$ws.Cells.Item(4,11).Value = 1.035918057762675
$ws.Cells.Item(4,12).Value = 1.033018855438276
$ws.Cells.Item(4,13).Value = 1.042824143620784
$ws.Cells.Item(4,14).Value = 1.037010353671228

$ws.Cells.Item(5,2).Value = 1.02
$ws.Cells.Item(5,3).Value = 1.031388255742771
$ws.Cells.Item(5,4).Value = 1.033727980926248
$ws.Cells.Item(5,5).Value = 1.030869378860774
$ws.Cells.Item(5,6).Value = 1.040993145082561
$ws.Cells.Item(5,9).Value = 1.033463529192469
$ws.Cells.Item(5,10).Value = 1.035844535765198
$ws.Cells.Item(5,11).Value = 1.036165802107071
$ws.Cells.Item(5,12).Value = 1.033314325954948
$ws.Cells.Item(5,13).Value = 1.043413045166071
$ws.Cells.Item(5,14).Value = 1.037315553819014

$ws.Cells.Item(6,2).Value = 1.02
$ws.Cells.Item(6,3).Value = 1.031455976376547
$ws.Cells.Item(6,4).Value = 1.033778392252951
$ws.Cells.Item(6,5).Value = 1.030927857077164
$ws.Cells.Item(6,6).Value = 1.041100840603253
$ws.Cells.Item(6,9).Value = 1.033483151479928
$ws.Cells.Item(6,10).Value = 1.035895673472698
$ws.Cells.Item(6,11).Value = 1.036207362394532
$ws.Cells.Item(6,12).Value = 1.033363905556787
$ws.Cells.Item(6,13).Value = 1.043511819020616
$ws.Cells.Item(6,14).Value = 1.037366764147924

$ws.Cells.Item(7,2).Value = 1.02
$ws.Cells.Item(7,3).Value = 1.030990147954533
$ws.Cells.Item(7,4).Value = 1.033431586246373
$ws.Cells.Item(7,5).Value = 1.030525612611129
$ws.Cells.Item(7,6).Value = 1.04035971832049
$ws.Cells.Item(7,9).Value = 1.033347981742479
$ws.Cells.Item(7,10).Value = 1.035543843025501
$ws.Cells.Item(7,11).Value = 1.035921370622992
$ws.Cells.Item(7,12).Value = 1.033022805624384
$ws.Cells.Item(7,13).Value = 1.042832019631965
$ws.Cells.Item(7,14).Value = 1.037014434061136

$ws.Cells.Item(8,2).Value = 1.02
$ws.Cells.Item(8,3).Value = 1.02903879111589
$ws.Cells.Item(8,4).Value = 1.031977753537116
$ws.Cells.Item(8,5).Value = 1.028840822704963
$ws.Cells.Item(8,6).Value = 1.037247016497547
$ws.Cells.Item(8,9).Value = 1.032776900169578
$ws.Cells.Item(8,10).Value = 1.03406826360571
$ws.Cells.Item(8,11).Value = 1.034720558433272
$ws.Cells.Item(8,12).Value = 1.031592496986029
$ws.Cells.Item(8,13).Value = 1.039975053977641
$ws.Cells.Item(8,14).Value = 1.035536759149317

$ws.Cells.Item(9,2).Value = 1.02
$ws.Cells.Item(9,3).Value = 1.02558346333761
$ws.Cells.Item(9,4).Value = 1.029399443955228
$ws.Cells.Item(9,5).Value = 1.025858371611212
$ws.Cells.Item(9,6).Value = 1.031702864424283
$ws.Cells.Item(9,9).Value = 1.031747258996574
$ws.Cells.Item(9,10).Value = 1.031448773559039
$ws.Cells.Item(9,11).Value = 1.032583705013107
$ws.Cells.Item(9,12).Value = 1.029054391526712
$ws.Cells.Item(9,13).Value = 1.034879524026838
$ws.Cells.Item(9,14).Value = 1.03291354912632

$ws.Cells.Item(10,2).Value = 1.02
$ws.Cells.Item(10,3).Value = 1.023267968470128
$ws.Cells.Item(10,4).Value = 1.027669032808473
$ws.Cells.Item(10,5).Value = 1.023860388001257
$ws.Cells.Item(10,6).Value = 1.027964225437221
$ws.Cells.Item(10,9).Value = 1.031044815980535
$ws.Cells.Item(10,10).Value = 1.029688939233127
$ws.Cells.Item(10,11).Value = 1.031144663425244
$ws.Cells.Item(10,12).Value = 1.027349928172122
$ws.Cells.Item(10,13).Value = 1.0314387827356
$ws.Cells.Item(10,14).Value = 1.031151215633811

$ws.Cells.Item(11,2).Value = 1.02
$ws.Cells.Item(11,3).Value = 1.022262353437898
$ws.Cells.Item(11,4).Value = 1.026916904067228
$ws.Cells.Item(11,5).Value = 1.022992826410458
$ws.Cells.Item(11,6).Value = 1.026334643107049
$ws.Cells.Item(11,9).Value = 1.030736775859835
$ws.Cells.Item(11,10).Value = 1.028923592548137
$ws.Cells.Item(11,11).Value = 1.030518010045564
$ws.Cells.Item(11,12).Value = 1.026608829121612
$ws.Cells.Item(11,13).Value = 1.029937951510898
$ws.Cells.Item(11,14).Value = 1.030384782068741

$ws.Cells.Item(12,2).Value = 1.02
$ws.Cells.Item(12,3).Value = 1.021888361116577
$ws.Cells.Item(12,4).Value = 1.026637092704146
$ws.Cells.Item(12,5).Value = 1.022670201484266
$ws.Cells.Item(12,6).Value = 1.025727678046926
$ws.Cells.Item(12,9).Value = 1.030621766808194
$ws.Cells.Item(12,10).Value = 1.028638798837636
$ws.Cells.Item(12,11).Value = 1.03028470266675
$ws.Cells.Item(12,12).Value = 1.026333083489104
$ws.Cells.Item(12,13).Value = 1.029378777565773
$ws.Cells.Item(12,14).Value = 1.030099583918505

$ws.Cells.Item(13,2).Value = 1.02
$ws.Cells.Item(13,3).Value = 1.021968604858996
$ws.Cells.Item(13,4).Value = 1.026697133127267
$ws.Cells.Item(13,5).Value = 1.022739422710852
$ws.Cells.Item(13,6).Value = 1.025857950364224
$ws.Cells.Item(13,9).Value = 1.030646463412505
$ws.Cells.Item(13,10).Value = 1.028699911291197
$ws.Cells.Item(13,11).Value = 1.030334772487919
$ws.Cells.Item(13,12).Value = 1.026392253214543
$ws.Cells.Item(13,13).Value = 1.029498799955264
$ws.Cells.Item(13,14).Value = 1.030160783158759

$ws.Cells.Item(14,2).Value = 1.02
$ws.Cells.Item(14,3).Value = 1.022231448617995
$ws.Cells.Item(14,4).Value = 1.026893783746284
$ws.Cells.Item(14,5).Value = 1.02296616582037
$ws.Cells.Item(14,6).Value = 1.026284505480467
$ws.Cells.Item(14,9).Value = 1.030727281234934
$ws.Cells.Item(14,10).Value = 1.028900061871176
$ws.Cells.Item(14,11).Value = 1.030498735860105
$ws.Cells.Item(14,12).Value = 1.026586045512894
$ws.Cells.Item(14,13).Value = 1.029891764961762
$ws.Cells.Item(14,14).Value = 1.03036121797552

$ws.Cells.Item(15,2).Value = 1.02
$ws.Cells.Item(15,3).Value = 1.02239333371008
$ws.Cells.Item(15,4).Value = 1.027014888464536
$ws.Cells.Item(15,5).Value = 1.023105819880058
$ws.Cells.Item(15,6).Value = 1.026547097819637
$ws.Cells.Item(15,9).Value = 1.030776997481943
$ws.Cells.Item(15,10).Value = 1.029023313382591
$ws.Cells.Item(15,11).Value = 1.030599687231732
$ws.Cells.Item(15,12).Value = 1.026705384985268
$ws.Cells.Item(15,13).Value = 1.030133657084022
$ws.Cells.Item(15,14).Value = 1.030484644518218

$ws.Cells.Item(16,2).Value = 1.02
$ws.Cells.Item(16,3).Value = 1.023334643730026
$ws.Cells.Item(16,4).Value = 1.027718888303726
$ws.Cells.Item(16,5).Value = 1.023917913327859
$ws.Cells.Item(16,6).Value = 1.028072144857946
$ws.Cells.Item(16,9).Value = 1.031065177325174
$ws.Cells.Item(16,10).Value = 1.029739661824588
$ws.Cells.Item(16,11).Value = 1.03118617704356
$ws.Cells.Item(16,12).Value = 1.027399047303376
$ws.Cells.Item(16,13).Value = 1.03153815258673
$ws.Cells.Item(16,14).Value = 1.031202010257168

$ws.Cells.Item(17,2).Value = 1.02
$ws.Cells.Item(17,3).Value = 1.023924292990033
$ws.Cells.Item(17,4).Value = 1.028159719525011
$ws.Cells.Item(17,5).Value = 1.024426662685754
$ws.Cells.Item(17,6).Value = 1.029025856856624
$ws.Cells.Item(17,9).Value = 1.031244902102438
$ws.Cells.Item(17,10).Value = 1.030188110538083
$ws.Cells.Item(17,11).Value = 1.031553113024406
$ws.Cells.Item(17,12).Value = 1.027833338767291
$ws.Cells.Item(17,13).Value = 1.032416183740772
$ws.Cells.Item(17,14).Value = 1.031651095819274

$ws.Cells.Item(18,2).Value = 1.02
$ws.Cells.Item(18,3).Value = 1.024267937766
$ws.Cells.Item(18,4).Value = 1.028416574630625
$ws.Cells.Item(18,5).Value = 1.024723174630017
$ws.Cells.Item(18,6).Value = 1.029581109633541
$ws.Cells.Item(18,9).Value = 1.031349358842659
$ws.Cells.Item(18,10).Value = 1.0304493627839
$ws.Cells.Item(18,11).Value = 1.031766799721117
$ws.Cells.Item(18,12).Value = 1.028086359495849
$ws.Cells.Item(18,13).Value = 1.032927269720299
$ws.Cells.Item(18,14).Value = 1.031912719073246

$ws.Cells.Item(19,2).Value = 1.02
$ws.Cells.Item(19,3).Value = 1.024385063347409
$ws.Cells.Item(19,4).Value = 1.028504109335565
$ws.Cells.Item(19,5).Value = 1.024824238332186
$ws.Cells.Item(19,6).Value = 1.029770263174806
$ws.Cells.Item(19,9).Value = 1.031384912708434
$ws.Cells.Item(19,10).Value = 1.030538389021668
$ws.Cells.Item(19,11).Value = 1.031839603810095
$ws.Cells.Item(19,12).Value = 1.028172583418856
$ws.Cells.Item(19,13).Value = 1.033101359448134
$ws.Cells.Item(19,14).Value = 1.032001871738484

$ws.Cells.Item(20,2).Value = 1.02
$ws.Cells.Item(20,3).Value = 1.023861059051642
$ws.Cells.Item(20,4).Value = 1.028112450944052
$ws.Cells.Item(20,5).Value = 1.024372102862453
$ws.Cells.Item(20,6).Value = 1.028923639644532
$ws.Cells.Item(20,9).Value = 1.031225658034319
$ws.Cells.Item(20,10).Value = 1.03014002943808
$ws.Cells.Item(20,11).Value = 1.031513779585315
$ws.Cells.Item(20,12).Value = 1.027786773895192
$ws.Cells.Item(20,13).Value = 1.032322088767561
$ws.Cells.Item(20,14).Value = 1.031602946438595

$ws.Cells.Item(21,2).Value = 1.02
$ws.Cells.Item(21,3).Value = 1.022154060542977
$ws.Cells.Item(21,4).Value = 1.026835887185353
$ws.Cells.Item(21,5).Value = 1.022899406030102
$ws.Cells.Item(21,6).Value = 1.026158942048845
$ws.Cells.Item(21,9).Value = 1.030703498720418
$ws.Cells.Item(21,10).Value = 1.028841136670517
$ws.Cells.Item(21,11).Value = 1.030450467722503
$ws.Cells.Item(21,12).Value = 1.026528991519102
$ws.Cells.Item(21,13).Value = 1.029776093857751
$ws.Cells.Item(21,14).Value = 1.030302209094317

$ws.Cells.Item(22,2).Value = 1.02
$ws.Cells.Item(22,3).Value = 1.021078122645745
$ws.Cells.Item(22,4).Value = 1.026030726379311
$ws.Cells.Item(22,5).Value = 1.021971294563415
$ws.Cells.Item(22,6).Value = 1.024411002604803
$ws.Cells.Item(22,9).Value = 1.030371783915347
$ws.Cells.Item(22,10).Value = 1.028021516171027
$ws.Cells.Item(22,11).Value = 1.029778790270916
$ws.Cells.Item(22,12).Value = 1.025735458752042
$ws.Cells.Item(22,13).Value = 1.028165473541917
$ws.Cells.Item(22,14).Value = 1.029481424639706

$ws.Cells.Item(23,2).Value = 1.02
$ws.Cells.Item(23,3).Value = 1.021648756062354
$ws.Cells.Item(23,4).Value = 1.02645780070146
$ws.Cells.Item(23,5).Value = 1.022463512976673
$ws.Cells.Item(23,6).Value = 1.02533855275921
$ws.Cells.Item(23,9).Value = 1.030547957931946
$ws.Cells.Item(23,10).Value = 1.028456296006474
$ws.Cells.Item(23,11).Value = 1.030135158937826
$ws.Cells.Item(23,12).Value = 1.026156386028954
$ws.Cells.Item(23,13).Value = 1.029020244901756
$ws.Cells.Item(23,14).Value = 1.02991682191239

$ws.Cells.Item(24,2).Value = 1.02
$ws.Cells.Item(24,3).Value = 1.023889632635355
$ws.Cells.Item(24,4).Value = 1.028133810428698
$ws.Cells.Item(24,5).Value = 1.024396756817944
$ws.Cells.Item(24,6).Value = 1.028969830387016
$ws.Cells.Item(24,9).Value = 1.031234354754882
$ws.Cells.Item(24,10).Value = 1.030161756207247
$ws.Cells.Item(24,11).Value = 1.031531553726547
$ws.Cells.Item(24,12).Value = 1.027807815465398
$ws.Cells.Item(24,13).Value = 1.032364609492501
$ws.Cells.Item(24,14).Value = 1.031624704062266

$ws.Cells.Item(25,2).Value = 1.02
$ws.Cells.Item(25,3).Value = 1.026478804851145
$ws.Cells.Item(25,4).Value = 1.030067997998146
$ws.Cells.Item(25,5).Value = 1.026631075234737
$ws.Cells.Item(25,6).Value = 1.033143453252651
$ws.Cells.Item(25,9).Value = 1.03201624386582
$ws.Cells.Item(25,10).Value = 1.032128317462427
$ws.Cells.Item(25,11).Value = 1.033138650971642
$ws.Cells.Item(25,12).Value = 1.029712701420169
$ws.Cells.Item(25,13).Value = 1.036204354778863
$ws.Cells.Item(25,14).Value = 1.033594058059996
